# revision on main board
# Adds two new BOM rows (35, 36) with a yellow highlight style, shifts the
# "Total" row down to row 37, and extends the SUM formula to cover the new
# rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: push the old "Total" row (36) down to row 37, freeing up
#     rows 35 and 36 for the two new BOM lines. -----------------------------
$ws.Range("A36").EntireRow.Insert()

# --- New row 35: J6 connector (not populated / not added to cart) --------
$ws.Range("A35").Value = 1
$ws.Range("E35").Value = "J6"
$ws.Range("G35").Value = "C29275"
$ws.Range("H35").Value = "C29275"
$ws.Range("I35").Value = "BOOMELE(Boom Precision Elec)"
$ws.Range("J35").Value = 0.0088769999999999995
$ws.Range("K35").Value = "KF2510 Header Male Pin 0.100""(2.54mm) 3 3P Wire To Board / Wire To Wire Connector RoHS"
$ws.Range("L35").Formula = "=A35*J35"
$ws.Range("M35").Value = "not added to cart"

# --- New row 36: F1 resettable fuse ---------------------------------------
$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "F1206"
$ws.Range("E36").Value = "F1"
$ws.Range("I36").Value = "SOCAY"
$ws.Range("H36").Value = "SCF075-1206R"
$ws.Range("G36").Value = "C183290"
$ws.Range("J36").Value = 0.054459
$ws.Range("K36").Value = "8V 1.5A SMD1206 PTC Resettable Fuses RoHS"
$ws.Range("L36").Formula = "=A36*J36"

# --- Apply the new yellow-highlight style to just the populated cells ----
foreach ($addr in @("A35","E35","G35","H35","I35","K35","M35",
                     "A36","B36","E36","G36","H36","I36","K36")) {
    $ws.Range($addr).Interior.Color = 65535
}
foreach ($addr in @("J35","J36")) {
    $ws.Range($addr).Interior.Color = 65535
    $ws.Range($addr).NumberFormat = """$""#,##0.000000"
}
foreach ($addr in @("L35","L36")) {
    $ws.Range($addr).Interior.Color = 65535
    $ws.Range($addr).NumberFormat = """$""#,##0.00"
}

# M26 also picks up the new yellow style per the diff
$ws.Range("M26").Interior.Color = 65535

# --- Total row (now row 37): extend the SUM to include the new rows ------
$ws.Range("L37").Formula = "=SUM(L2:L36)"

$ws.Range("L37").Select()
